{"js": "// Insert a new \"Bedroom\" scene-heading paragraph at the very start of the\n// document body. The paragraph carries the same body formatting (Calibri,\n// 12pt, space-after 200) as the rest of the script, with the run itself\n// bolded, matching a scene-slug style heading.\n//\n// We build the new paragraph (plus the required paragraph-end marker) as a\n// Flat OPC OOXML fragment and insert it with Body.insertOoxml so that the\n// exact paragraph/run structure (including the placeholder trailing run)\n// is produced, then hand off to Word's own parser to merge it in.\nconst body = context.document.body;\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr>' +\n  '<w:spacing w:after=\"200\" w:lineRule=\"auto\"/>' +\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:eastAsia=\"Calibri\" w:hAnsi=\"Calibri\"/>' +\n  '<w:sz w:val=\"24\"/>' +\n  '<w:szCs w:val=\"24\"/>' +\n  '</w:rPr>' +\n  '</w:pPr>' +\n  '<w:r>' +\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:eastAsia=\"Calibri\" w:hAnsi=\"Calibri\"/>' +\n  '<w:b w:val=\"1\"/>' +\n  '<w:sz w:val=\"24\"/>' +\n  '<w:szCs w:val=\"24\"/>' +\n  '<w:rtl w:val=\"0\"/>' +\n  '</w:rPr>' +\n  '<w:t xml:space=\"preserve\">Bedroom</w:t>' +\n  '</w:r>' +\n  '<w:r>' +\n  '<w:rPr>' +\n  '<w:rtl w:val=\"0\"/>' +\n  '</w:rPr>' +\n  '</w:r>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nbody.insertOoxml(flatOpcXml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Insert a new \"Bedroom\" scene-heading paragraph at the very start of the\n# document body. The paragraph carries the same body formatting (Calibri,\n# 12pt, space-after 200) as the rest of the script, with the run itself\n# bolded, matching a scene-slug style heading.\n#\n# A collapsed Range at the very start of the story is fed a Flat OPC OOXML\n# fragment via Range.InsertXML so the exact paragraph/run structure\n# (including the placeholder trailing run) is produced without disturbing\n# any of the existing paragraphs that follow it.\n$d = $word.ActiveDocument\n\n$rng = $d.Range(0, 0)\n\n$xml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:after=\"200\" w:lineRule=\"auto\"/>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:eastAsia=\"Calibri\" w:hAnsi=\"Calibri\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:eastAsia=\"Calibri\" w:hAnsi=\"Calibri\"/>\n                <w:b w:val=\"1\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:rtl w:val=\"0\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">Bedroom</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rtl w:val=\"0\"/>\n              </w:rPr>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$rng.InsertXML($xml)\n"}
